$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 306, shifting existing rows 306..407 down to 307..408
$ws.Rows.Item(306).Insert()

# Populate the new row 306 with its data
$ws.Cells.Item(306, 1).Value = 10
$ws.Cells.Item(306, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(306, 3).Value = "La Araucanía"
$ws.Cells.Item(306, 4).Value = 44985
$ws.Cells.Item(306, 4).NumberFormat = $ws.Cells.Item(305, 4).NumberFormat
$ws.Cells.Item(306, 5).Value = 9
$ws.Cells.Item(306, 6).Value = 100112001
$ws.Cells.Item(306, 7).Value = "Berenjena"
$ws.Cells.Item(306, 8).Value = "Sin especificar"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 80
$ws.Cells.Item(306, 11).Value = 14000
$ws.Cells.Item(306, 12).Value = 14000
$ws.Cells.Item(306, 13).Value = 14000
$ws.Cells.Item(306, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(306, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(306, 16).Value = 350
$ws.Cells.Item(306, 17).Value = 40
$ws.Cells.Item(306, 18).Value = "Hortaliza"
